# Update the "deal_value" column (D2:D81) on the active sheet with
# recalculated sensitivity values following an update to probabilities.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    1879.608760016147,
    1872.70833147886,
    1852.61774475008,
    1822.988234314474,
    1786.674403157196,
    1745.874615795742,
    1702.25949194861,
    1657.082776385526,
    1611.273388107098,
    1565.509850870335,
    1520.279361933564,
    1475.924038266588,
    1432.676762002742,
    1390.688753272352,
    1350.050653083261,
    1310.808564687997,
    1272.976206577506,
    1236.544082360581,
    1201.486371240702,
    1167.766082327189,
    1135.338890051725,
    1104.155970052505,
    1074.166079308563,
    1045.317066244124,
    1017.556952066964,
    990.834690646679,
    965.100688341877,
    940.3071454584576,
    916.4082660068872,
    893.3603710061341,
    871.1219419005878,
    849.6536140605947,
    828.9181353297315,
    808.8803007833022,
    789.5068719840052,
    770.7664868425579,
    752.6295645456497,
    735.0682087736591,
    718.0561114986716,
    701.5684589557166,
    685.5818408602438,
    670.0741635597434,
    655.0245675242919,
    640.4133493747892,
    626.2218884996455,
    612.4325782062731,
    599.0287612819266,
    585.994669790823,
    573.3153689048364,
    560.9767045483815,
    548.9652546310896,
    537.2682836414482,
    525.8737003789743,
    514.7700186100285,
    503.9463204420151,
    493.3922222215153,
    483.0978427734255,
    473.0537738098259,
    463.2510523488398,
    453.6811349950307,
    1132.513751629532,
    1204.613673277342,
    1230.871634042808,
    1222.189455187303,
    1187.837790446574,
    1135.338890051725,
    1070.58347751069,
    998.0579537941217,
    921.1014240532265,
    842.1484953439607,
    762.9385077627255,
    684.6862534358118,
    608.2164709370491,
    534.0673353313099,
    462.5688003560213,
    393.9012125093449,
    328.1387813240341,
    265.2815927118892,
    205.2790404393421,
    148.0468738621572
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 4).Value = $newValues[$i]
}

$wb.Save()
